$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data as described in the commit.
# For columns D (Price) and E (Volume(1h)) we prefix values with a leading
# apostrophe so Excel stores them as literal text (matching the original
# inline-string cell contents) instead of auto-converting to numbers, and we
# reset the cell style to "Normal" afterwards so no extra text-format style is
# introduced (the source cells carry no explicit style).

# Row 2
$ws.Range("D2").Value = "'89.880.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.90%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'3.052.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.70%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.06%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'212.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.44%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'613.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.45%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'0.363"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -8.53%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("E8").Value = "'  +14.86%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.04%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'3.047.95"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.73%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("E11").Value = "'  +21.20%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'0.188"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +4.82%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'0.0000240"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -4.54%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'5.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.33%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'89.490.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.68%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'32.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.15%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("E17").Value = "'  -3.27%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'3.047.28"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -3.27%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'3.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.57%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'0.0000219"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -4.68%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'13.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.06%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'425.65"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.67%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'8.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.45%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'5.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +2.22%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'5.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.60%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'84.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.86%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'11.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.97%  "
$ws.Range("E27").Style = "Normal"

# Row 29
$ws.Range("E29").Value = "'  +3.10%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("E30").Value = "'  +0.24%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("D31").Value = "'8.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.36%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'3.75"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -7.22%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'504.67"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.94%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'6.68"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -5.36%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "'22.89"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +4.05%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("E36").Value = "'  -2.18%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("E37").Value = "'  -3.91%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("D38").Value = "'0.132"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -10.82%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'22.25"
$ws.Range("D39").Style = "Normal"

# Row 40
$ws.Range("E40").Value = "'  -0.06%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("E41").Value = "'  -0.02%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'0.364"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.01%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("E43").Value = "'  +4.76%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'1.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.57%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'146.69"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.93%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "'0.0691"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +11.97%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").Value = "'43.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.93%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'4.22"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +7.09%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'162.14"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.68%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'1.22"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +2.61%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'23.88"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.80%  "
$ws.Range("E51").Style = "Normal"
